$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rebuild the breakpoint scale as a regular
# 0,12.5,20,25,30,35,40,45,50,55 ladder (was 12.5,20,25,30,35,42.5,50,max)
# and push "Ti le" / "Bat cap" headers two columns further right.
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 12.5
$ws.Range("E1").Value = 20
$ws.Range("F1").Value = 25
$ws.Range("G1").Value = 30
$ws.Range("H1").Value = 35
$ws.Range("I1").Value = 40
$ws.Range("J1").Value = 45
$ws.Range("K1").Value = 50
$ws.Range("L1").Value = 55
$ws.Range("M1").Value = "Ti le"
$ws.Range("N1").Value = "Bat cap"

# --- Row 2 (Hai): two extra ratio columns inserted before the trailing pair
$ws.Range("I2").Value = 0.018
$ws.Range("J2").Value = 0.019
$ws.Range("K2").Value = 0.02
$ws.Range("L2").Value = 0.021

# --- Row 3 (Tien)
$ws.Range("I3").Value = 0.017
$ws.Range("J3").Value = 0.018
$ws.Range("K3").Value = 0.019
$ws.Range("L3").Value = 0.02

# --- Row 4 (Minh): flat 0.01 row extended two more columns
$ws.Range("K4").Value = 0.01
$ws.Range("L4").Value = 0.01

# --- Row 5 (Cuong): ratios shift right, plus the "0.7 / *" footnote pair
$ws.Range("I5").Value = 0.017
$ws.Range("J5").Value = 0.018
$ws.Range("K5").Value = 0.019
$ws.Range("L5").Value = 0.02
$ws.Range("M5").Value = 0.7
$ws.Range("N5").Value = "*"

# --- Row 6 (Duc)
$ws.Range("I6").Value = 0.017
$ws.Range("J6").Value = 0.018
$ws.Range("K6").Value = 0.019
$ws.Range("L6").Value = 0.02

# --- Column widths: C keeps manual width (no more bestFit), D:L become a
# narrow fixed-width block sized for the short numeric headers.
$ws.Columns.Item(3).ColumnWidth = 11.75
$ws.Range("D1:L1").EntireColumn.ColumnWidth = 5.17

# --- Selection moves to D1 (matches the saved cursor position)
$ws.Range("D1").Select() | Out-Null
